# Append two new data rows (117 and 118) to the active sheet, matching the
# source dataset's next two days of COVID figures (columns A-R; S is left
# blank, same as the trailing column in every other existing row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(116, 1873, 67, 986, 47, 2973, 37, 30, 67, 21036, 14010, 4053, 288, 16983, 2, 449, 1053, 1211),
    @(117, 2095, 77, 945, 49, 3166, 42, 35, 77, 22031, 14712, 4153, 193, 17878, 2, 222, 1022, 895)
)

$startRow = 117
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $rowValues[$col]
    }
}
